# "cambios de las fracciones"
# Update reporting period dates on the "Reporte de Formatos" sheet:
#   B8 (Fecha de inicio) : 2022-07-01 -> 2022-10-01
#   C8 (Fecha de término): 2022-09-30 -> 2022-12-31
#   S8 (Fecha de validación)   : 2022-10-10 -> 2023-01-10
#   T8 (Fecha de actualización): 2022-10-10 -> 2023-01-10
# Also move the active selection/view to match where Excel left the
# cursor after the edit (merged cell G3:I3, scrolled to column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

$ws.Range("B8").Value = 44835
$ws.Range("C8").Value = 44926
$ws.Range("S8").Value = 44936
$ws.Range("T8").Value = 44936

$ws.Activate()
$ws.Range("G3:I3").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 6
